# Auto-generated edit script applying the Sheets diff.
# For each touched cell: set new value, or clear the cell if it was removed,
# or set the value on a newly-added cell.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 24391788
$ws.Range("I137").Value = 1010.125
$ws.Range("K137").Value = 3030.375
$ws.Range("M137").Value = -480.375
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2708.2917
$ws.Range("I3").Value = 1833.1666
$ws.Range("K3").Value = 1833.1666
$ws.Range("M3").Value = -1718.1666
$ws.Range("H61").Value = 1905.8485
$ws.Range("I61").Value = 1236.2593
$ws.Range("J61").Value = 4919
$ws.Range("K61").Value = 1236.2593
$ws.Range("L61").Value = 4919
$ws.Range("M61").Value = -1024.2593
$ws.Range("N61").Value = -5343
$ws.Range("H74").Value = 3750.4187
$ws.Range("I74").Value = 862.6667
$ws.Range("J74").Value = 5829.6
$ws.Range("K74").Value = 862.6667
$ws.Range("L74").Value = 5829.6
$ws.Range("M74").Value = 11.33330000000001
$ws.Range("N74").Value = -7577.6
$ws.Range("H77").Value = 3750.4187
$ws.Range("I77").Value = 862.6667
$ws.Range("J77").Value = 5829.6
$ws.Range("K77").Value = 4313.3335
$ws.Range("L77").Value = 29148
$ws.Range("M77").Value = 54.66650000000027
$ws.Range("N77").Value = -37884
$ws.Range("H132").Value = 1423.8572
$ws.Range("I132").Value = 1330.5667
$ws.Range("J132").Value = 1657.0834
$ws.Range("K132").Value = 3991.7001
$ws.Range("L132").Value = 4971.2502
$ws.Range("M132").Value = -1461.7001
$ws.Range("N132").Value = -10031.2502
$ws.Range("H136").Value = 1905.8485
$ws.Range("I136").Value = 1236.2593
$ws.Range("J136").Value = 4919
$ws.Range("K136").Value = 3708.7779
$ws.Range("L136").Value = 14757
$ws.Range("M136").Value = -1158.7779
$ws.Range("N136").Value = -19857
$ws.Range("H139").Value = 46442.4
$ws.Range("J139").Value = 46442.4
$ws.Range("L139").Value = 46442.4
$ws.Range("N139").Value = -56722.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 50000252
$ws.Range("I7").Value = 50000252
$ws.Range("K7").Value = 50000252
$ws.Range("M7").Value = -50000139
$ws.Range("H105").Value = 7166.6665
$ws.Range("I105").Value = 10000
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 10000
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -8253
$ws.Range("N105").Value = -4994
$ws.Range("H134").Value = 36505.594
$ws.Range("I134").Value = 39729.344
$ws.Range("K134").Value = 119188.032
$ws.Range("M134").Value = -116653.032

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1071.3611
$ws.Range("I31").Value = 989.34485
$ws.Range("J31").Value = 1411.1428
$ws.Range("K31").Value = 989.34485
$ws.Range("L31").Value = 1411.1428
$ws.Range("M31").Value = -694.34485
$ws.Range("N31").Value = -2001.1428
$ws.Range("H34").Value = 1071.3611
$ws.Range("I34").Value = 989.34485
$ws.Range("J34").Value = 1411.1428
$ws.Range("K34").Value = 989.34485
$ws.Range("L34").Value = 1411.1428
$ws.Range("M34").Value = -787.34485
$ws.Range("N34").Value = -1815.1428
$ws.Range("H58").Value = 2138.5264
$ws.Range("I58").Value = 1480.4348
$ws.Range("J58").Value = 3147.6
$ws.Range("K58").Value = 1480.4348
$ws.Range("L58").Value = 3147.6
$ws.Range("M58").Value = -1277.4348
$ws.Range("N58").Value = -3553.6
$ws.Range("H132").Value = 1702.1316
$ws.Range("J132").Value = 4000.75
$ws.Range("L132").Value = 12002.25
$ws.Range("N132").Value = -17062.25
$ws.Range("H134").Value = 1528.6765
$ws.Range("I134").Value = 1413.6207
$ws.Range("J134").Value = 2196
$ws.Range("K134").Value = 4240.8621
$ws.Range("L134").Value = 6588
$ws.Range("M134").Value = -1705.8621
$ws.Range("N134").Value = -11658
$ws.Range("H136").Value = 2138.5264
$ws.Range("I136").Value = 1480.4348
$ws.Range("J136").Value = 3147.6
$ws.Range("K136").Value = 4441.3044
$ws.Range("L136").Value = 9442.799999999999
$ws.Range("M136").Value = -1891.3044
$ws.Range("N136").Value = -14542.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 37037652
$ws.Range("J34").Value = 38462172
$ws.Range("L34").Value = 115386516
$ws.Range("N34").Value = -115386684
$ws.Range("H40").Value = 479.5909
$ws.Range("I40").Value = 71.75
$ws.Range("J40").Value = 969
$ws.Range("K40").Value = 287
$ws.Range("L40").Value = 3876
$ws.Range("M40").Value = -218
$ws.Range("N40").Value = -4014
$ws.Range("H68").Value = 956.8823
$ws.Range("I68").Value = 864.8333
$ws.Range("J68").Value = 1007.0909
$ws.Range("K68").Value = 2594.4999
$ws.Range("L68").Value = 3021.2727
$ws.Range("M68").Value = -1783.4999
$ws.Range("N68").Value = -4643.2727
$ws.Range("H70").Value = 3558655
$ws.Range("I70").Value = 5926125
$ws.Range("K70").Value = 17778375
$ws.Range("M70").Value = -17778060
$ws.Range("H71").Value = 956.8823
$ws.Range("I71").Value = 864.8333
$ws.Range("J71").Value = 1007.0909
$ws.Range("K71").Value = 7783.4997
$ws.Range("L71").Value = 9063.8181
$ws.Range("M71").Value = -3727.4997
$ws.Range("N71").Value = -17175.8181
$ws.Range("H73").Value = 3558655
$ws.Range("I73").Value = 5926125
$ws.Range("K73").Value = 17778375
$ws.Range("M73").Value = -17777283
$ws.Range("H94").Value = 2745.2
$ws.Range("I94").Value = 1224
$ws.Range("J94").Value = 5027
$ws.Range("K94").Value = 3672
$ws.Range("L94").Value = 15081
$ws.Range("M94").Value = -2996
$ws.Range("N94").Value = -16433
$ws.Range("H107").Value = 64906.84
$ws.Range("I107").Value = 200290.6
$ws.Range("J107").Value = 38871.5
$ws.Range("K107").Value = 600871.8
$ws.Range("L107").Value = 116614.5
$ws.Range("M107").Value = -598951.8
$ws.Range("N107").Value = -120454.5
$ws.Range("H133").Value = 7452.143
$ws.Range("J133").Value = 8053.2354
$ws.Range("L133").Value = 24159.7062
$ws.Range("N133").Value = -34279.7062

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5967.5557
$ws.Range("I70").Value = 5986.857
$ws.Range("K70").Value = 5986.857
$ws.Range("M70").Value = -5716.857
$ws.Range("H73").Value = 5967.5557
$ws.Range("I73").Value = 5986.857
$ws.Range("K73").Value = 5986.857
$ws.Range("M73").Value = -5050.857
$ws.Range("H105").Value = 16250
$ws.Range("J105").Value = 16250
$ws.Range("L105").Value = 16250
$ws.Range("N105").Value = -23238
$ws.Range("H132").Value = 2453.0605
$ws.Range("I132").Value = 2376.85
$ws.Range("J132").Value = 2570.3076
$ws.Range("K132").Value = 7130.549999999999
$ws.Range("L132").Value = 7710.9228
$ws.Range("M132").Value = -4600.549999999999
$ws.Range("N132").Value = -12770.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4379.325
$ws.Range("I132").Value = 4791.4814
$ws.Range("J132").Value = 3523.3076
$ws.Range("K132").Value = 14374.4442
$ws.Range("L132").Value = 10569.9228
$ws.Range("M132").Value = -11844.4442
$ws.Range("N132").Value = -15629.9228
$ws.Range("H136").Value = 1624.1666
$ws.Range("I136").Value = 549
$ws.Range("K136").Value = 1647
$ws.Range("M136").Value = 903

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3960.111
$ws.Range("I132").Value = 4923.5835
$ws.Range("K132").Value = 14770.7505
$ws.Range("M132").Value = -12240.7505
$ws.Range("H136").Value = 10470.857
$ws.Range("I136").Value = 11871.556
$ws.Range("K136").Value = 35614.66800000001
$ws.Range("M136").Value = -33064.66800000001
